$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2€")

# "Obv: Without mint symbol" varieties now in stock -> flip the 0 -> 1
# quantity markers for the 2021 (I16/I17), 2022 (H19) and 2023 (I20/I21) rows.
$ws.Range("I16").Value = 1
$ws.Range("I17").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("I20").Value = 1
$ws.Range("I21").Value = 1

# Move the active sheet's selection (bottom-right frozen pane) to L24,
# matching where the author left the cursor after the edit.
$ws.Activate()
$ws.Range("L24").Select()
